# Insert a new "is_normal_for_donor" column between "Sample_ID" (D) and
# "relative_file_path" (previously E, now shifts to F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing column E ("relative_file_path" data) one column to the
# right, opening up a blank column E for the new field.
$ws.Columns.Item(5).Insert()

# Header for the new column.
$ws.Range("E1").Value = "is_normal_for_donor"

# Only the "test2" sample row (row 3) carries a value for the new column.
$ws.Range("E3").Value = "Y"

# Match the column width Excel would have used for the newly inserted column.
$ws.Columns.Item(5).ColumnWidth = 17.5

# Reflect the new selection left behind after the edit.
$ws.Range("E5").Select()
